$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The quiz-results block that lived under the "UNIT 1" Quiz/Results columns
# (A/C, rows 5-8) actually belongs under the "UNIT 2" Quiz/Results columns
# (F/H, rows 5-8) -- move (not copy) the values there, carrying the same
# cell formatting along, and remove the now-empty source cells entirely.

$styleA = $ws.Range("A5").Style
$styleC = $ws.Range("C5").Style

$ws.Range("F5:F8").Value2 = $ws.Range("A5:A8").Value2
$ws.Range("F5:F8").Style = $styleA

$ws.Range("H5:H8").Value2 = $ws.Range("C5:C8").Value2
$ws.Range("H5:H8").Style = $styleC

$ws.Range("A5:A8").Clear()
$ws.Range("C5:C8").Clear()
